$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82

# Copy formatting from the row above (row 81) down into the new row,
# then set the new values - this mirrors the "daily update" pattern of
# extending the table by one day.
$ws.Range("A81:D81").Copy() | Out-Null
$ws.Range("A82:D82").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 46031
$ws.Cells.Item($row, 2).Value = 184
$ws.Cells.Item($row, 3).Value = 195
$ws.Cells.Item($row, 4).Value = 187
